$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values include figures like "42.402.39" (multi-dot) and plain
# decimals like "236.93" that Excel would otherwise auto-convert to numbers.
# Force the whole Price column to Text first so values are stored as strings,
# matching the original inline-string cell contents, then restore default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.402.39'
$ws.Range('D3').Value = '2.245.50'
$ws.Range('D5').Value = '236.93'
$ws.Range('D6').Value = '0.633'
$ws.Range('D7').Value = '69.74'
$ws.Range('D10').Value = '0.0995'
$ws.Range('D12').Value = '36.64'
$ws.Range('D14').Value = '6.76'
$ws.Range('D15').Value = '2.579.25'
$ws.Range('D16').Value = '15.07'
$ws.Range('D17').Value = '0.870'
$ws.Range('D18').Value = '2.251.48'
$ws.Range('D19').Value = '42.255.58'
$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('D22').Value = '73.42'
$ws.Range('D23').Value = '236.41'
$ws.Range('D24').Value = '2.02'
$ws.Range('D29').Value = '10.04'
$ws.Range('D30').Value = '170.95'
$ws.Range('D31').Value = '20.60'
$ws.Range('D33').Value = '0.127'
$ws.Range('D34').Value = '0.0721'
$ws.Range('D36').Value = '4.71'
$ws.Range('D37').Value = '3.76'
$ws.Range('D38').Value = '23.01'
$ws.Range('D42').Value = '65.53'
$ws.Range('D43').Value = '9.36'
$ws.Range('D46').Value = '0.193'
$ws.Range('D49').Value = '10.24'
$ws.Range('D51').Value = '2.34'

$ws.Range("D2:D51").Style = "Normal"

$ws.Range('E2').Value = '  -3.02%  '
$ws.Range('E3').Value = '  -3.91%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('E6').Value = '  -4.72%  '
$ws.Range('E7').Value = '  -2.77%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -6.44%  '
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  +1.69%  '
$ws.Range('E12').Value = '  +14.00%  '
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('E14').Value = '  -4.83%  '
$ws.Range('E15').Value = '  -3.95%  '
$ws.Range('E16').Value = '  -5.92%  '
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('E18').Value = '  -3.79%  '
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('E21').Value = '  -4.55%  '
$ws.Range('E22').Value = '  -5.40%  '
$ws.Range('E23').Value = '  -5.69%  '
$ws.Range('E24').Value = '  +6.66%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('E27').Value = '  -2.26%  '
$ws.Range('E28').Value = '  -3.90%  '
$ws.Range('E29').Value = '  -2.30%  '
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('E31').Value = '  -6.80%  '
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('E33').Value = '  -4.90%  '
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  -6.54%  '
$ws.Range('E37').Value = '  +1.44%  '
$ws.Range('E38').Value = '  +23.48%  '
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('E41').Value = '  -6.27%  '
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('E43').Value = '  +2.87%  '
$ws.Range('E44').Value = '  -16.06%  '
$ws.Range('E45').Value = '  -3.07%  '
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('E47').Value = '  +13.70%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  +10.86%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('E51').Value = '  -2.41%  '
